# Update marksheet "Corr/total" marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking row, Right column -> 5 (was 3)
$ws.Range("B11").Value = 5

# B12: Total row, Right column -> 70 (was 42)
$ws.Range("B12").Value = 70

# E12: Total row, Max column, "corr/total" text -> 70/140 (was 40/84)
$ws.Range("E12").Value = "70/140"

$wb.Save()
